$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) "this question should be answered with " + "c#" -> single run
#    (drops the spellStart/spellEnd proofErr pair)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(2)
$p.Range.InsertXML('<w:p ' + $wNs + '>' +
  '<w:r><w:rPr><w:i/><w:color w:val="0000FF"/></w:rPr>' +
  '<w:t>this question should be answered with c#</w:t></w:r>' +
  '</w:p>')

# ---------------------------------------------------------------------------
# 2) "Write the program in the most efficient way you can." + " Can you tell
#    what is your solution complexity? Is it O(n)?" -> single run
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(6)
$p.Range.InsertXML('<w:p ' + $wNs + '>' +
  '<w:pPr><w:spacing w:before="220" w:after="220"/></w:pPr>' +
  '<w:r><w:t>Write the program in the most efficient way you can. Can you tell what is your solution complexity? Is it O(n)?</w:t></w:r>' +
  '</w:p>')

# ---------------------------------------------------------------------------
# 3) "number 5 repeated" + " 2 times," -> single run
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(15)
$p.Range.InsertXML('<w:p ' + $wNs + '>' +
  '<w:r><w:t>number 5 repeated 2 times,</w:t></w:r>' +
  '</w:p>')

# ---------------------------------------------------------------------------
# 4) "that would " + "be,  X" + "=[1,1,5,5,3,3,3]." -> single run
#    (drops the gramStart/gramEnd proofErr pair)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(17)
$p.Range.InsertXML('<w:p ' + $wNs + '>' +
  '<w:r><w:t>that would be,  X=[1,1,5,5,3,3,3].</w:t></w:r>' +
  '</w:p>')

# ---------------------------------------------------------------------------
# 5) "explana" + "tion:" -> "explanation:" (second example's heading)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(22)
$p.Range.InsertXML('<w:p ' + $wNs + '>' +
  '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' +
  '<w:r><w:t>explanation:</w:t></w:r>' +
  '</w:p>')

# ---------------------------------------------------------------------------
# 6) code-sample table cell: "[" + "]" -> "[]" (drops gramStart/gramEnd)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(25)
$cons = '<w:rFonts w:ascii="Consolas" w:eastAsia="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>'
$p.Range.InsertXML('<w:p ' + $wNs + '>' +
  '<w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F5F5F5"/><w:spacing w:line="325" w:lineRule="auto"/>' +
  '<w:rPr>' + $cons + '<w:color w:val="777777"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr>' + $cons + '<w:color w:val="4B69C6"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>public</w:t></w:r>' +
  '<w:r><w:rPr>' + $cons + '<w:color w:val="333333"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr>' + $cons + '<w:color w:val="4B69C6"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>int</w:t></w:r>' +
  '<w:r><w:rPr>' + $cons + '<w:color w:val="333333"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr>' + $cons + '<w:b/><w:color w:val="AA3731"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>Challenge</w:t></w:r>' +
  '<w:r><w:rPr>' + $cons + '<w:color w:val="777777"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>(</w:t></w:r>' +
  '<w:r><w:rPr>' + $cons + '<w:color w:val="4B69C6"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>int</w:t></w:r>' +
  '<w:r><w:rPr>' + $cons + '<w:color w:val="777777"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>[]</w:t></w:r>' +
  '<w:r><w:rPr>' + $cons + '<w:color w:val="333333"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> input</w:t></w:r>' +
  '<w:r><w:rPr>' + $cons + '<w:color w:val="777777"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>){</w:t></w:r>' +
  '</w:p>')

# ---------------------------------------------------------------------------
# 7) Add the closing answer paragraph text before the _GoBack bookmark.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs($d.Paragraphs.Count)
$p.Range.InsertXML('<w:p ' + $wNs + '>' +
  '<w:r><w:t xml:space="preserve">Please find the answer in the attached code. Required method is implemented in the </w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Solver</w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t>class.</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>')

# ---------------------------------------------------------------------------
# 8) Footer: "All right reserved for " + "Axceligent" + " Solutions" ->
#    single run (drops spellStart/spellEnd proofErr pair)
# ---------------------------------------------------------------------------
$footer = $d.Sections(1).Footers.Item(1)
$fp = $footer.Range.Paragraphs(1)
$fp.Range.InsertXML('<w:p ' + $wNs + '>' +
  '<w:pPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr>' +
  '<w:tabs><w:tab w:val="center" w:pos="4680"/><w:tab w:val="right" w:pos="9360"/></w:tabs>' +
  '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:color w:val="000000"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>All right reserved for Axceligent Solutions</w:t></w:r>' +
  '</w:p>')
